$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 24
$ws.Range("B3").Value = 840000
$ws.Range("B4").Value = 4285714.285714285
$ws.Range("B5").Value = 3452000
$ws.Range("B6").Value = 670000
$ws.Range("B11").Value = 100000
$ws.Range("B35").Value = 9347714.285714285
$ws.Range("B38").Value = 9347714.285714285

$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 21
$ws.Range("B3").Value = 2250000
$ws.Range("B13").Value = 1500000
$ws.Range("B23").Value = 2250000
$ws.Range("B31").Value = -2700000
$ws.Range("B32").Value = 1500000
$ws.Range("B33").Value = 2920000
$ws.Range("B34").Value = 1720000.000000001

$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 22
$ws.Range("B3").Value = 2553571.428571429
$ws.Range("B8").Value = 900000
$ws.Range("B10").Value = 150000
$ws.Range("B34").Value = 803571.4285714286
$ws.Range("B37").Value = 903571.4285714286

$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = 26
$ws.Range("B3").Value = 1857142.857142857
$ws.Range("B13").Value = 2785714.285714286
$ws.Range("B23").Value = 2785714.285714286
$ws.Range("B31").Value = 1857142.857142857
$ws.Range("B32").Value = 2785714.285714286
$ws.Range("B33").Value = 2785714.285714286
$ws.Range("B34").Value = 7428571.428571429

$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = 26
$ws.Range("B3").Value = 5571428.571428572
$ws.Range("B31").Value = 5571428.571428572
$ws.Range("B34").Value = 5571428.571428572

$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = 24

$ws = $wb.Worksheets.Item(8)
$ws.Range("B2").Value = 26
$ws.Range("B3").Value = 6500000
$ws.Range("B31").Value = 0
$ws.Range("B34").Value = 0

$ws = $wb.Worksheets.Item(9)
$ws.Range("B2").Value = 25
$ws.Range("B3").Value = 875000
$ws.Range("B4").Value = 7142857.142857144
$ws.Range("B7").Value = 5586000
$ws.Range("B32").Value = 14003857.14285714
$ws.Range("B35").Value = 22338142.85714286

$ws = $wb.Worksheets.Item(10)
$ws.Range("B2").Value = 26
$ws.Range("B3").Value = 910000
$ws.Range("B4").Value = 3095238.095238095
$ws.Range("B5").Value = 4950000
$ws.Range("B7").Value = 770000
$ws.Range("B8").Value = 1620000
$ws.Range("B15").Value = 3095238.095238095
$ws.Range("B26").Value = 3095238.095238095
$ws.Range("B35").Value = 13360238.0952381
$ws.Range("B36").Value = 14745238.0952381
$ws.Range("B37").Value = 11695238.0952381
$ws.Range("B38").Value = 39800714.28571428

$ws = $wb.Worksheets.Item(11)
$ws.Range("B2").Value = 15.5
$ws.Range("B3").Value = 542500
$ws.Range("B4").Value = 1660714.285714286
$ws.Range("B8").Value = 1620000
$ws.Range("B13").Value = -600000
$ws.Range("B16").Value = 1660714.285714286
$ws.Range("B27").Value = 1660714.285714286
$ws.Range("B36").Value = 4913214.285714284
$ws.Range("B37").Value = 4082142.857142857
$ws.Range("B38").Value = 4710714.285714285
$ws.Range("B39").Value = 13706071.42857143

$ws = $wb.Worksheets.Item(12)
$ws.Range("B2").Value = 26
$ws.Range("B3").Value = 910000
$ws.Range("B4").Value = 13928571.42857143
$ws.Range("B15").Value = 9285714.285714285
$ws.Range("B26").Value = 13928571.42857143
$ws.Range("B35").Value = 10426571.42857143
$ws.Range("B36").Value = 9285714.285714285
$ws.Range("B37").Value = 13928571.42857143
$ws.Range("B38").Value = 33640857.14285715

$ws = $wb.Worksheets.Item(13)
$ws.Range("B12").Value = 23.5
$ws.Range("B13").Value = 6714285.714285715
$ws.Range("B32").Value = -535714.2857142854
$ws.Range("B34").Value = -535714.2857142854

$ws = $wb.Worksheets.Item(14)
$ws.Range("B12").Value = 24
$ws.Range("B13").Value = 4285714.285714285
$ws.Range("B32").Value = 3285714.285714285
$ws.Range("B34").Value = 3985714.285714285

$ws = $wb.Worksheets.Item(15)
$ws.Range("B12").Value = 22
$ws.Range("B13").Value = 3142857.142857143
$ws.Range("B32").Value = 792857.1428571432
$ws.Range("B34").Value = 792857.1428571432

$ws = $wb.Worksheets.Item(16)
$ws.Range("B13").Value = 22
$ws.Range("B14").Value = 770000
$ws.Range("B15").Value = 3928571.428571429
$ws.Range("B36").Value = 5898571.428571429
$ws.Range("B38").Value = 5898571.428571429

$ws = $wb.Worksheets.Item(17)
$ws.Range("B12").Value = 20
$ws.Range("B13").Value = 700000
$ws.Range("B14").Value = 2142857.142857143
$ws.Range("B33").Value = 1512857.142857143
$ws.Range("B35").Value = 1512857.142857143

$ws = $wb.Worksheets.Item(18)
$ws.Range("B5").Value = 1900000
$ws.Range("B24").Value = 25
$ws.Range("B25").Value = 875000
$ws.Range("B26").Value = 8035714.285714285
$ws.Range("B35").Value = 2350000
$ws.Range("B37").Value = 4630714.285714285
$ws.Range("B38").Value = 12680714.28571429

$ws = $wb.Worksheets.Item(19)
$ws.Range("B22").Value = 20.5
$ws.Range("B23").Value = 717500
$ws.Range("B24").Value = 4026785.714285714
$ws.Range("B34").Value = 2894285.714285715
$ws.Range("B35").Value = 2894285.714285715

$ws = $wb.Worksheets.Item(20)
$ws.Range("B22").Value = 22
$ws.Range("B23").Value = 770000
$ws.Range("B24").Value = 3142857.142857143
$ws.Range("B34").Value = -1537142.857142857
$ws.Range("B35").Value = -1537142.857142857

$ws = $wb.Worksheets.Item(21)
$ws.Range("B24").Value = 25
$ws.Range("B25").Value = 4464285.714285715
$ws.Range("B36").Value = 4514285.714285715
$ws.Range("B37").Value = 4514285.714285715

$ws = $wb.Worksheets.Item(22)
$ws.Range("B22").Value = 20.5
$ws.Range("B23").Value = 717500
$ws.Range("B24").Value = 2196428.571428571
$ws.Range("B34").Value = 2913928.571428571
$ws.Range("B35").Value = 2913928.571428571
